$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 6661.25
$ws.Range("J38").Value = 12000
$ws.Range("L38").Value = 36000
$ws.Range("N38").Value = -36744
$ws.Range("H74").Value = 13999.75
$ws.Range("I74").Value = 7999.5
$ws.Range("K74").Value = 7999.5
$ws.Range("M74").Value = -7063.5
$ws.Range("H77").Value = 13999.75
$ws.Range("I77").Value = 7999.5
$ws.Range("K77").Value = 39997.5
$ws.Range("M77").Value = -35317.5
$ws.Range("H98").Value = 52631900
$ws.Range("I98").Value = 55555840
$ws.Range("J98").Value = 1000
$ws.Range("K98").Value = 55555840
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = -55554342
$ws.Range("N98").Value = -3996
$ws.Range("H122").Value = 52631900
$ws.Range("I122").Value = 55555840
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 166667520
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -166665070
$ws.Range("N122").Value = -7900
$ws.Range("H132").Value = 2735.7273
$ws.Range("I132").Value = 2115.8064
$ws.Range("K132").Value = 6347.4192
$ws.Range("M132").Value = -3817.4192
$ws.Range("H135").Value = 2386.8
$ws.Range("I135").Value = 2402.5
$ws.Range("K135").Value = 21622.5
$ws.Range("M135").Value = -19087.5
$ws.Range("H137").Value = 6040.091
$ws.Range("I137").Value = 4085.5
$ws.Range("J137").Value = 6773.0625
$ws.Range("K137").Value = 12256.5
$ws.Range("L137").Value = 20319.1875
$ws.Range("M137").Value = -9706.5
$ws.Range("N137").Value = -25419.1875
$ws.Range("H141").Value = 4056.5454
$ws.Range("I141").Value = 4537.0625
$ws.Range("J141").Value = 2775.1667
$ws.Range("K141").Value = 13611.1875
$ws.Range("L141").Value = 8325.500100000001
$ws.Range("M141").Value = -8431.1875
$ws.Range("N141").Value = -18685.5001

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 33335
$ws.Range("J10").Value = 33335
$ws.Range("L10").Value = 33335
$ws.Range("N10").Value = -33675
$ws.Range("H63").Value = 9825
$ws.Range("I63").Value = 9750
$ws.Range("K63").Value = 9750
$ws.Range("M63").Value = -9064
$ws.Range("H66").Value = 9825
$ws.Range("I66").Value = 9750
$ws.Range("K66").Value = 48750
$ws.Range("M66").Value = -45318
$ws.Range("H102").Value = 2576.5
$ws.Range("I102").Value = 2576.5
$ws.Range("K102").Value = 2576.5
$ws.Range("M102").Value = -954.5
$ws.Range("H110").Value = 2333.3333
$ws.Range("I110").Value = 2333.3333
$ws.Range("K110").Value = 2333.3333
$ws.Range("M110").Value = -288.3332999999998
$ws.Range("H132").Value = 6409.6904
$ws.Range("I132").Value = 2610.5925
$ws.Range("K132").Value = 7831.7775
$ws.Range("M132").Value = -5301.7775

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 1588
$ws.Range("I80").Value = 465
$ws.Range("J80").Value = 2037.2
$ws.Range("K80").Value = 465
$ws.Range("L80").Value = 2037.2
$ws.Range("M80").Value = 533
$ws.Range("N80").Value = -4033.2
$ws.Range("H83").Value = 1588
$ws.Range("I83").Value = 465
$ws.Range("J83").Value = 2037.2
$ws.Range("K83").Value = 2325
$ws.Range("L83").Value = 10186
$ws.Range("M83").Value = 2667
$ws.Range("N83").Value = -20170
$ws.Range("H86").Value = 2200.4167
$ws.Range("I86").Value = 2205.5
$ws.Range("K86").Value = 2205.5
$ws.Range("M86").Value = -1082.5
$ws.Range("H89").Value = 2200.4167
$ws.Range("I89").Value = 2205.5
$ws.Range("K89").Value = 11027.5
$ws.Range("M89").Value = -5411.5
$ws.Range("H107").Value = 3201.25
$ws.Range("I107").Value = 3285
$ws.Range("K107").Value = 3285
$ws.Range("M107").Value = -1365
$ws.Range("H134").Value = 51230.855
$ws.Range("I134").Value = 3791.7
$ws.Range("J134").Value = 1000014
$ws.Range("K134").Value = 11375.1
$ws.Range("L134").Value = 3000042
$ws.Range("M134").Value = -8840.099999999999
$ws.Range("N134").Value = -3005112

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 694524
$ws.Range("I31").Value = 13572.23
$ws.Range("J31").Value = 1116065.5
$ws.Range("K31").Value = 13572.23
$ws.Range("L31").Value = 1116065.5
$ws.Range("M31").Value = -13277.23
$ws.Range("N31").Value = -1116655.5
$ws.Range("H34").Value = 694524
$ws.Range("I34").Value = 13572.23
$ws.Range("J34").Value = 1116065.5
$ws.Range("K34").Value = 13572.23
$ws.Range("L34").Value = 1116065.5
$ws.Range("M34").Value = -13370.23
$ws.Range("N34").Value = -1116469.5
$ws.Range("H58").Value = 1450.4286
$ws.Range("I58").Value = 1215
$ws.Range("J58").Value = 1764.3334
$ws.Range("K58").Value = 1215
$ws.Range("L58").Value = 1764.3334
$ws.Range("M58").Value = -1012
$ws.Range("N58").Value = -2170.3334
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()
$ws.Range("H94").Value = 9837.333000000001
$ws.Range("I94").Value = 11000
$ws.Range("J94").Value = 9604.799999999999
$ws.Range("K94").Value = 11000
$ws.Range("L94").Value = 9604.799999999999
$ws.Range("M94").Value = -10549
$ws.Range("N94").Value = -10506.8
$ws.Range("H132").Value = 4465.591
$ws.Range("I132").Value = 4418.9473
$ws.Range("K132").Value = 13256.8419
$ws.Range("M132").Value = -10726.8419
$ws.Range("H136").Value = 1450.4286
$ws.Range("I136").Value = 1215
$ws.Range("J136").Value = 1764.3334
$ws.Range("K136").Value = 3645
$ws.Range("L136").Value = 5293.0002
$ws.Range("M136").Value = -1095
$ws.Range("N136").Value = -10393.0002

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 4232606
$ws.Range("I4").Value = 6050859
$ws.Range("K4").Value = 18152577
$ws.Range("M4").Value = -18152465
$ws.Range("H5").Value = 271
$ws.Range("I5").Value = 271
$ws.Range("K5").Value = 813
$ws.Range("M5").Value = -701
$ws.Range("H12").Value = 634191.4
$ws.Range("I12").Value = 1739.3334
$ws.Range("J12").Value = 1055826.1
$ws.Range("K12").Value = 5218.0002
$ws.Range("L12").Value = 3167478.3
$ws.Range("M12").Value = -5045.0002
$ws.Range("N12").Value = -3167824.3
$ws.Range("H14").Value = 408.64285
$ws.Range("I14").Value = 408.64285
$ws.Range("K14").Value = 1225.92855
$ws.Range("M14").Value = -1052.92855
$ws.Range("H75").Value = 1549999
$ws.Range("J75").Value = 1549999
$ws.Range("L75").Value = 4649997
$ws.Range("N75").Value = -4651993
$ws.Range("H78").Value = 1549999
$ws.Range("J78").Value = 1549999
$ws.Range("L78").Value = 13949991
$ws.Range("N78").Value = -13959975
$ws.Range("H122").Value = 638.3043
$ws.Range("I122").Value = 596
$ws.Range("J122").Value = 665.5
$ws.Range("K122").Value = 5364
$ws.Range("L122").Value = 5989.5
$ws.Range("M122").Value = -2914
$ws.Range("N122").Value = -10889.5
$ws.Range("H131").Value = 1815.125
$ws.Range("I131").Value = 1980.5
$ws.Range("J131").Value = 1649.75
$ws.Range("K131").Value = 5941.5
$ws.Range("L131").Value = 4949.25
$ws.Range("M131").Value = -901.5
$ws.Range("N131").Value = -15029.25
$ws.Range("H135").Value = 271
$ws.Range("I135").Value = 271
$ws.Range("K135").Value = 2439
$ws.Range("M135").Value = 96

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H101").Value = 132813
$ws.Range("J101").Value = 132813
$ws.Range("L101").Value = 132813
$ws.Range("N101").Value = -139303
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()
$ws.Range("H133").Value = 94000
$ws.Range("J133").Value = 94000
$ws.Range("L133").Value = 94000
$ws.Range("N133").Value = -104120

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4125.4375
$ws.Range("I40").Value = 3250.5
$ws.Range("K40").Value = 3250.5
$ws.Range("M40").Value = -3114.5
$ws.Range("H93").Value = 41668000
$ws.Range("I93").Value = 50001196
$ws.Range("K93").Value = 50001196
$ws.Range("M93").Value = -49999948
$ws.Range("H122").Value = 7643.778
$ws.Range("I122").Value = 6558.8
$ws.Range("K122").Value = 19676.4
$ws.Range("M122").Value = -17226.4
$ws.Range("H128").Value = 111207.5
$ws.Range("J128").Value = 111207.5
$ws.Range("L128").Value = 111207.5
$ws.Range("N128").Value = -121167.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H74").Value = 7686.875
$ws.Range("J74").Value = 7686.875
$ws.Range("L74").Value = 7686.875
$ws.Range("N74").Value = -9558.875
$ws.Range("H75").Value = 87937
$ws.Range("J75").Value = 87937
$ws.Range("L75").Value = 87937
$ws.Range("N75").Value = -89809
$ws.Range("H77").Value = 7686.875
$ws.Range("J77").Value = 7686.875
$ws.Range("L77").Value = 23060.625
$ws.Range("N77").Value = -32420.625
$ws.Range("H78").Value = 87937
$ws.Range("J78").Value = 87937
$ws.Range("L78").Value = 263811
$ws.Range("N78").Value = -273171
